$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.528.24"
$ws.Range("E2").Value = "  +6.64%  "
$ws.Range("D3").Value = "2.579.52"
$ws.Range("E3").Value = "  +8.59%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "504.85"
$ws.Range("E5").Value = "  +5.61%  "
$ws.Range("D6").Value = "156.92"
$ws.Range("E6").Value = "  +6.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.620"
$ws.Range("E7").Value = "  +23.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "2.579.22"
$ws.Range("E9").Value = "  +8.60%  "
$ws.Range("D10").Value = "6.14"
$ws.Range("E10").Value = "  +13.23%  "
$ws.Range("E11").Value = "  +5.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.340"
$ws.Range("E12").Value = "  +5.35%  "
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("D14").Value = "3.012.50"
$ws.Range("E14").Value = "  +8.11%  "
$ws.Range("D15").Value = "59.408.73"
$ws.Range("E15").Value = "  +6.27%  "
$ws.Range("D16").Value = "21.93"
$ws.Range("E16").Value = "  +7.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000137"
$ws.Range("E17").Value = "  +3.94%  "
$ws.Range("D18").Value = "2.571.55"
$ws.Range("E18").Value = "  +8.29%  "
$ws.Range("D19").Value = "4.72"
$ws.Range("E19").Value = "  +2.79%  "
$ws.Range("D20").Value = "334.75"
$ws.Range("E20").Value = "  +6.14%  "
$ws.Range("D21").Value = "10.37"
$ws.Range("E21").Value = "  +6.99%  "
$ws.Range("E22").Value = "  +6.91%  "
$ws.Range("D23").Value = "1.01"
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("D24").Value = "60.43"
$ws.Range("E24").Value = "  +6.52%  "
$ws.Range("D25").Value = "0.414"
$ws.Range("E25").Value = "  +4.84%  "
$ws.Range("E26").Value = "  +7.44%  "
$ws.Range("D27").Value = "2.666.44"
$ws.Range("E27").Value = "  +7.31%  "
$ws.Range("D28").Value = "0.995"
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("D29").Value = "7.45"
$ws.Range("E29").Value = "  +2.94%  "
$ws.Range("D30").Value = "0.0₃0821"
$ws.Range("E30").Value = "  +6.90%  "
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").Value = "156.34"
$ws.Range("E32").Value = "  +6.26%  "
$ws.Range("D33").Value = "19.42"
$ws.Range("E33").Value = "  +7.55%  "
$ws.Range("D34").Value = "1.56"
$ws.Range("E34").Value = "  +5.46%  "
$ws.Range("D35").Value = "5.48"
$ws.Range("E35").Value = "  +8.27%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "1.19"
$ws.Range("E36").Value = "  +7.79%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.90"
$ws.Range("E37").Value = "  +8.93%  "
$ws.Range("D38").Value = "0.861"
$ws.Range("E38").Value = "  +2.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.80"
$ws.Range("E39").Value = "  +12.48%  "
$ws.Range("E40").Value = "  +7.24%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "34.98"
$ws.Range("E41").Value = "  +4.71%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "291.39"
$ws.Range("E42").Value = "  +14.73%  "
$ws.Range("E43").Value = "  +7.26%  "
$ws.Range("D44").Value = "0.624"
$ws.Range("E44").Value = "  +7.48%  "
$ws.Range("E45").Value = "  +4.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("D47").Value = "19.41"
$ws.Range("E47").Value = "  +15.20%  "
$ws.Range("E48").Value = "  +6.69%  "
$ws.Range("D49").Value = "4.78"
$ws.Range("E49").Value = "  +5.29%  "
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("D51").Value = "0.714"
$ws.Range("E51").Value = "  +12.51%  "
